# 2.5.0 add variable name
# Update the "st_level" sheet headers to include type-prefixed variable
# names, and flag the sheet's JSON export as a generated item class.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("st_level")

# Row 1 JSON config cell: enable IsGenItemClass for st_level.
$ws.Range("C1").Value = '{"IsStringId":false,"IsGenItemClass":true,"JSONName":"st_levelJSON"}'

# Row 2 header cells: prefix field names with their variable type.
$ws.Range("B2").Value = "float:row"
$ws.Range("C2").Value = "string:col"
$ws.Range("H2").Value = "bool:testfloat"

# Match the author's final selection in the sheet.
$ws.Range("F4").Select()
